# The author renamed the "questiontags" sheet to "tags" and retitled its
# header cell B1 from "questiontag" to "tag". They also cleared the
# (invisible) number-format override that had been lingering on the
# A358:C416 block of the "warm" sheet, leaving those cells on the default
# "General" style. Everything else in the produced diff (the large
# cascade of shared-string index shifts across every sheet) is simply
# Excel's own shared-string table compaction falling out of that single
# text edit, so it does not need to be reproduced by hand here.

$wb = $excel.ActiveWorkbook

# 1) Rename "questiontags" -> "tags"
$wsTags = $wb.Worksheets.Item("questiontags")
$wsTags.Name = "tags"

# 2) Retitle the header cell on the renamed sheet
$wsTags.Range("B1").Value = "tag"

# 3) Drop the stray number-format style from the old questiontag rows on
#    "warm" (A358:C416) so they fall back to the default/General style
$wsWarm = $wb.Worksheets.Item("warm")
$wsWarm.Range("A358:C416").ClearFormats()

# 4) Leave the workbook with "warm" scrolled down near the bottom of its
#    data and "tags" as the active/selected sheet, with its header cell
#    selected - matching the saved view state of the edited workbook.
$wsWarm.Activate()
$wsWarm.Range("B315").Select()
$excel.ActiveWindow.ScrollRow = 313

$wsTags.Activate()
$wsTags.Range("B1").Select()
